# "Generate Report for Handoff"
#
# The localization status report is regenerated: the display text of each
# handoff file now shows the new handoff GUID/commit
# (9d5b49ee-7b3f-4506-88d0-827263ea55ea) instead of the old one
# (f3e9ec3e-7508-4815-a86a-99a1826266e1), and the new target-file hash
# (5881faf911e4e46acb2a3d06ffabd1dc628117a1) instead of the old one
# (ccbd6a1a4e98ccb3afe4dc1d1cfe1012b8e0afc8); the handoff timestamps are
# also refreshed. The underlying hyperlink addresses (the actual commit
# the links resolve to) are left untouched. This touches the "Overview",
# "zh-cn" and "de-de" sheets.

$wb = $excel.ActiveWorkbook

$newGuid = "9d5b49ee-7b3f-4506-88d0-827263ea55ea"
$newHash = "5881faf911e4e46acb2a3d06ffabd1dc628117a1"

# Original hyperlink target addresses (unchanged by this edit - only the
# display text / cell text is refreshed).
$mdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/1bc4a3c8a32be0fb3028bc9a96fa2495e614df2f/e2e/f3e9ec3e-7508-4815-a86a-99a1826266e1.md"
$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/47b8c5b3b13213fbfd5d3ebd8320aae221757aa7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/f3e9ec3e-7508-4815-a86a-99a1826266e1.ccbd6a1a4e98ccb3afe4dc1d1cfe1012b8e0afc8.zh-cn.xlf"
$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4a89851f7980fbe7a27784f003548657477bde3d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/f3e9ec3e-7508-4815-a86a-99a1826266e1.ccbd6a1a4e98ccb3afe4dc1d1cfe1012b8e0afc8.de-de.xlf"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $mdUrl, "", "", "$newGuid.md")

$ws.Range("D2").Value = "2016-42-14 08:42:08"

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $mdUrl, "", "", "$newGuid.md")
$ws.Hyperlinks.Add($ws.Range("B2"), $mdUrl, "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D2"), $zhXlfUrl, "", "", "$newGuid.$newHash.zh-cn.xlf")

$ws.Range("E2").Value = "2016-03-14 08:42:05"

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $mdUrl, "", "", "$newGuid.md")
$ws.Hyperlinks.Add($ws.Range("B2"), $mdUrl, "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D2"), $deXlfUrl, "", "", "$newGuid.$newHash.de-de.xlf")

$ws.Range("E2").Value = "2016-03-14 08:42:08"
